$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout used spacer columns (C, E, G, I, K) between each "Day" column
# (B, D, F, H, J, L). Delete the spacer columns (rightmost first so indices
# stay stable) to compact the six day-columns into contiguous B:G, which also
# naturally shifts every cell value / merged range onto its new reference and
# shrinks the used range to A1:G30.
$ws.Columns.Item(11).Delete()
$ws.Columns.Item(9).Delete()
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(5).Delete()
$ws.Columns.Item(3).Delete()

# Give the six remaining data columns a uniform width (matches width="35" in
# the saved workbook once Excel quantizes the character width).
$ws.Range("B1:G1").EntireColumn.ColumnWidth = 34.14

# Row 1: replace the generic "Day N" headers with actual dates.
$ws.Range("B1").Value = "14 Jul (Monday)"
$ws.Range("C1").Value = "15 Jul (Tuesday)"
$ws.Range("D1").Value = "16 Jul (Wednesday)"
$ws.Range("E1").Value = "17 Jul (Thursday)"
$ws.Range("F1").Value = "18 Jul (Friday)"
$ws.Range("G1").Value = "19 Jul (Saturday)"

# Small copy fixes: lowercase "lesson" -> "Lesson" for consistency.
$ws.Range("D7").Value = "Private Lesson with Stephane RETY " + [char]10 + "(Room Stephane)"
$ws.Range("B19").Value = "Private Lesson with Ivy CHUANG " + [char]10 + "(Room Ivy)"
$ws.Range("E19").Value = "Private Lesson with Ivy CHUANG " + [char]10 + "(Room Ivy)"

# Writing into those wrapped multi-line cells nudges Excel into stamping an
# explicit row height; auto-fitting the rows back restores the implicit
# (default) height so the saved XML has no stray ht/customHeight attributes.
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(19).AutoFit()
